# Adds default s3-bucket values ("open" / "scratch") to the sample
# data-transfer manifest sheet, and moves the active selection to E5
# (matching the author's final cursor position after editing column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default bucket values for the two data rows (column E = "s3-bucket")
$ws.Range("E3").Value = "open"
$ws.Range("E4").Value = "scratch"

# Leave the selection where the author left it after making the edit
$ws.Range("E5").Select()
